$wb = $excel.ActiveWorkbook

# --- Data edit: updated Harvest_Cost_Inflation_Rate row (row 6) on the
#     inflation_rates sheet with newer projected values ---
$wsInflation = $wb.Worksheets.Item("inflation_rates")
$wsInflation.Range("B6").Value = -0.0121
$wsInflation.Range("C6").Value = 0.0031
$wsInflation.Range("D6").Value = 0.0156
$wsInflation.Range("E6").Value = 0.0182
$wsInflation.Range("F6").Value = 0.0147

# --- Leave-behind UI state: per-sheet selections, reflecting where the
#     author's cursor ended up on each tab ---
$wb.Worksheets.Item("Assumption (1)").Range("A28").Select() | Out-Null
$wb.Worksheets.Item("Base_cost").Range("I12").Select() | Out-Null
$wb.Worksheets.Item("Alt1_cost").Range("L19").Select() | Out-Null
$wb.Worksheets.Item("Insurance_aph_price").Range("C4").Select() | Out-Null
$wb.Worksheets.Item("Incurance_alternatives_frac").Range("E12").Select() | Out-Null
$wb.Worksheets.Item("Incurance_alternatives_prem").Range("C11").Select() | Out-Null
$wb.Worksheets.Item("planted_acres").Range("E4").Select() | Out-Null
$wb.Worksheets.Item("Alt2_cost").Range("B2:H5").Select() | Out-Null
$wb.Worksheets.Item("Alt3_cost").Range("C10").Select() | Out-Null

# inflation_rates is the sheet left active/selected when the file was saved
$wsInflation.Range("H18").Select() | Out-Null
$wsInflation.Activate() | Out-Null
